$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = '[Guang%Chen%xref no email%0,    Di%Wu%xref no email%0,    Wei%Guo%xref no email%0,    Yong%Cao%xref no email%0,    Da%Huang%xref no email%0,    Hongwu%Wang%xref no email%0,    Tao%Wang%xref no email%0,    Xiaoyun%Zhang%xref no email%0,    Huilong%Chen%xref no email%0,    Haijing%Yu%xref no email%0,    Xiaoping%Zhang%xref no email%0,    Minxia%Zhang%xref no email%0,    Shiji%Wu%xref no email%0,    Jianxin%Song%xref no email%0,    Tao%Chen%xref no email%0,    Meifang%Han%xref no email%0,    Shusheng%Li%xref no email%0,    Xiaoping%Luo%xref no email%0,    Jianping%Zhao%xref no email%0,    Qin%Ning%xref no email%0]'
$ws.Cells.Item(3, 5).Value2 = '[Chai%X%coreGivesNoEmail%2,          Deliang%Huang%coreGivesNoEmail%2,          Hong%Yu%coreGivesNoEmail%2,          Jun%Chen%coreGivesNoEmail%0,          Lei%Liu%coreGivesNoEmail%0,          National%Health Commission of the People\u2019s Republic of China%coreGivesNoEmail%2,          Pengcheng%Ou%coreGivesNoEmail%2,          Qing%He%coreGivesNoEmail%2,          Qingxian%Cai%coreGivesNoEmail%2,          Rossaint%J%coreGivesNoEmail%2,          Yang%Fu%coreGivesNoEmail%2,          Yiming%Zhang%coreGivesNoEmail%2,          Yinan%Su%coreGivesNoEmail%2,          Zhang%Xia%coreGivesNoEmail%2,          Zhenghua%Ma%coreGivesNoEmail%2,          Zhibin%Zhu%coreGivesNoEmail%2,          Zhiwei%Li%coreGivesNoEmail%2]'
$ws.Cells.Item(4, 5).Value2 = '[Yong%Gao%NULL%0,           Tuantuan%Li%NULL%0,           Mingfeng%Han%NULL%0,           Xiuyong%Li%NULL%0,           Dong%Wu%NULL%0,           Yuanhong%Xu%NULL%0,           Yulin%Zhu%NULL%0,           Yan%Liu%NULL%0,           Xiaowu%Wang%wangxiaowu19880218@126.com%0,           Linding%Wang%wanglinding@ahmu.edu.cn%0]'
$ws.Cells.Item(5, 5).Value2 = '[Zhongliang%Wang%NULL%0,           Bohan%Yang%NULL%0,           Qianwen%Li%NULL%0,           Lu%Wen%NULL%0,           Ruiguang%Zhang%zrg27@163.com%0]'
$ws.Cells.Item(6, 5).Value2 = '[Chuan%Qin%NULL%0,           Luoqi%Zhou%NULL%1,           Ziwei%Hu%NULL%1,           Shuoqi%Zhang%NULL%1,           Sheng%Yang%NULL%1,           Yu%Tao%NULL%0,           Cuihong%Xie%NULL%1,           Ke%Ma%NULL%0,           Ke%Shang%NULL%1,           Wei%Wang%NULL%0,           Dai-Shi%Tian%tiands@tjh.tjmu.edu.cn%1]'
$ws.Cells.Item(7, 5).Value2 = '[Chaomin%Wu%xref no email%0,    Xiaoyan%Chen%xref no email%0,    Yanping%Cai%xref no email%0,    Jia\u2019an%Xia%xref no email%0,    Xing%Zhou%xref no email%0,    Sha%Xu%xref no email%0,    Hanping%Huang%xref no email%0,    Li%Zhang%xref no email%0,    Xia%Zhou%xref no email%0,    Chunling%Du%xref no email%0,    Yuye%Zhang%xref no email%0,    Juan%Song%xref no email%0,    Sijiao%Wang%xref no email%0,    Yencheng%Chao%xref no email%0,    Zeyong%Yang%xref no email%0,    Jie%Xu%xref no email%0,    Xin%Zhou%xref no email%0,    Dechang%Chen%xref no email%0,    Weining%Xiong%xref no email%0,    Lei%Xu%xref no email%0,    Feng%Zhou%xref no email%0,    Jinjun%Jiang%xref no email%0,    Chunxue%Bai%xref no email%0,    Junhua%Zheng%xref no email%0,    Yuanlin%Song%xref no email%0]'
$ws.Cells.Item(9, 5).Value2 = '[Ruirui%Wang%NULL%0,           Min%Pan%NULL%0,           Xiumei%Zhang%NULL%0,           Mingfeng%Han%fyhmf@163.com%0,           Xiaoyun%Fan%13956988552@126.com%0,           Fengde%Zhao%NULL%0,           Manli%Miao%NULL%0,           Jing%Xu%NULL%0,           Minglong%Guan%NULL%0,           Xia%Deng%NULL%0,           Xu%Chen%NULL%0,           Leilei%Shen%NULL%0]'
$ws.Cells.Item(10, 5).Value2 = '[Zhe%Zhu%NULL%1,           Ting%Cai%NULL%1,           Lingyan%Fan%NULL%1,           Kehong%Lou%NULL%1,           Xin%Hua%NULL%1,           Zuoan%Huang%NULL%1,           Guosheng%Gao%NULL%1]'
$ws.Cells.Item(11, 5).Value2 = '[Xiaohua%Chen%NULL%1,           Binghong%Zhao%NULL%1,           Yueming%Qu%NULL%1,           Yurou%Chen%NULL%1,           Jie%Xiong%NULL%1,           Yong%Feng%NULL%2,           Dong%Men%NULL%1,           Qianchuan%Huang%NULL%1,           Ying%Liu%NULL%0,           Bo%Yang%NULL%0,           Jinya%Ding%NULL%1,           Feng%Li%fli222@whu.edu.cn%0]'
$ws.Cells.Item(12, 5).Value2 = '[Ming%Ding%NULL%1,           Qiang%Zhang%NULL%1,           Qing%Li%NULL%1,           Ting%Wu%NULL%1,           Ying-zi%Huang%NULL%1]'
$ws.Cells.Item(13, 5).Value2 = '[Peder L.%Myhre%p.l.myhre@medisin.uio.no%1,           Christian%Prebensen%christian.prebensen@gmail.com%1,           Heidi%Strand%Heidi.Strand@ahus.no%1,           Ragnhild%Røysland%ragnhiro@medisin.uio.no%1,           Christine M.%Jonassen%chrjon@so-hf.no%1,           Anbjørg%Rangberg%Anbjorg.Rangberg@so-hf.no%1,           Vibecke%Sørensen%Vibecke.Sorensen@ahus.no%1,           Signe%Søvik%signe.sovik@medisin.uio.no%1,           Helge%Røsjø%helge.rosjo@medisin.uio.no%1,           My%Svensson%m.h.s.svensson@medisin.uio.no%1,           Jan%Erik Berdal%jan-erik.berdal@hotmail.com%1,           Torbjørn%Omland%NULL%1]'
$ws.Cells.Item(14, 5).Value2 = '[Ana Carolina%Gadotti%NULL%1,           Marina%de Castro Deus%NULL%1,           Joao Paulo%Telles%NULL%1,           Rafael%Wind%NULL%1,           Marina%Goes%NULL%1,           Roberta%Garcia Charello Ossoski%NULL%1,           Alessandra Michalski%de Padua%NULL%1,           Lucia%de Noronha%NULL%1,           Andrea%Moreno-Amaral%NULL%1,           Cristina Pellegrino%Baena%NULL%1,           Felipe Francisco%Tuon%NULL%1]'
$ws.Cells.Item(15, 5).Value2 = '[Rocio%Laguna-Goya%NULL%1,           Alberto%Utrero-Rico%NULL%1,           Paloma%Talayero%NULL%1,           Maria%Lasa-Lazaro%NULL%1,           Angel%Ramirez-Fernandez%NULL%1,           Laura%Naranjo%NULL%1,           Alejandro%Segura-Tudela%NULL%1,           Oscar%Cabrera-Marante%NULL%1,           Edgar%Rodriguez de Frias%NULL%1,           Rocio%Garcia-Garcia%NULL%1,           Mario%Fernández-Ruiz%NULL%0,           Jose Maria%Aguado%NULL%1,           Joaquin%Martinez-Lopez%NULL%1,           Elena Ana%Lopez%NULL%1,           Mercedes%Catalan%NULL%1,           Antonio%Serrano%NULL%1,           Estela%Paz-Artal%NULL%1]'
$ws.Cells.Item(16, 5).Value2 = '[Tao%Chen%NULL%0,           Di%Wu%NULL%0,           Huilong%Chen%NULL%0,           Weiming%Yan%NULL%0,           Danlei%Yang%NULL%0,           Guang%Chen%NULL%0,           Ke%Ma%NULL%0,           Dong%Xu%NULL%0,           Haijing%Yu%NULL%0,           Hongwu%Wang%NULL%0,           Tao%Wang%NULL%0,           Wei%Guo%NULL%0,           Jia%Chen%NULL%0,           Chen%Ding%NULL%0,           Xiaoping%Zhang%NULL%0,           Jiaquan%Huang%NULL%0,           Meifang%Han%NULL%0,           Shusheng%Li%NULL%0,           Xiaoping%Luo%NULL%0,           Jianping%Zhao%NULL%0,           Qin%Ning%NULL%0]'
$ws.Cells.Item(17, 5).Value2 = '[Lang%Wang%NULL%0,           Wenbo%He%NULL%0,           Xiaomei%Yu%NULL%0,           Dalong%Hu%NULL%0,           Mingwei%Bao%NULL%0,           Huafen%Liu%NULL%0,           Jiali%Zhou%NULL%0,           Hong%Jiang%NULL%0]'
$ws.Cells.Item(18, 5).Value2 = '[Fei%Zhou%NULL%0,           Ting%Yu%NULL%0,           Ronghui%Du%NULL%0,           Guohui%Fan%NULL%0,           Ying%Liu%NULL%0,           Zhibo%Liu%NULL%0,           Jie%Xiang%NULL%0,           Yeming%Wang%NULL%0,           Bin%Song%NULL%0,           Xiaoying%Gu%NULL%0,           Lulu%Guan%NULL%0,           Yuan%Wei%NULL%0,           Hui%Li%NULL%0,           Xudong%Wu%NULL%0,           Jiuyang%Xu%NULL%0,           Shengjin%Tu%NULL%0,           Yi%Zhang%NULL%0,           Hua%Chen%NULL%0,           Bin%Cao%NULL%0]'
$ws.Cells.Item(19, 5).Value2 = '[Haiying%Sun%NULL%0,           Ruoqi%Ning%NULL%1,           Yu%Tao%NULL%0,           Chong%Yu%NULL%2,           Xiaoyan%Deng%NULL%1,           Caili%Zhao%NULL%1,           Silu%Meng%NULL%1,           Fangxu%Tang%89650793@qq.com%1,           Dong%Xu%89650793@qq.com%0,           Dong%Xu%89650793@qq.com%0]'
$ws.Cells.Item(20, 5).Value2 = '[Junli%Fan%NULL%1,           Hui%Wang%NULL%0,           Guangming%Ye%NULL%2,           Xiaoling%Cao%NULL%1,           Xianqun%Xu%NULL%1,           Wenbin%Tan%wenbin.tan@uscmed.sc.edu%1,           Yongxi%Zhang%znact1936@126.com%0]'
$ws.Cells.Item(21, 5).Value2 = '[Mathilda%Mandel%Mati.Mandel@sheba.health.gov.il%1,           Gil%Harari%NULL%1,           Michael%Gurevich%NULL%1,           Anat%Achiron%NULL%1]'
$ws.Cells.Item(22, 5).Value2 = '[Jose J.%Guirao%NULL%1,           Carmen M.%Cabrera%NULL%1,           Natalia%Jiménez%NULL%1,           Laura%Rincón%NULL%1,           José M.%Urra%NULL%1]'
$ws.Cells.Item(23, 5).Value2 = '[Chunjin%Ke%NULL%1,           Chong%Yu%NULL%0,           Daoyuan%Yue%NULL%1,           Xing%Zeng%NULL%1,           Zhiquan%Hu%NULL%1,           Chunguang%Yang%NULL%1]'
$ws.Cells.Item(24, 5).Value2 = '[Fuyang%Chen%NULL%0,           Wenwu%Sun%NULL%1,           Shengrong%Sun%NULL%1,           Zhiyu%Li%lizhiyu@whu.edu.cn%1,           Zhong%Wang%zhongwangchn@whu.edu.cn%1,           Li%Yu%yuliwhzxyy@163.com%2,           Li%Yu%yuliwhzxyy@163.com%0]'
$ws.Cells.Item(25, 5).Value2 = '[Hui%Chen%NULL%1,           Jun%Wang%NULL%0,           Nan%Su%NULL%1,           Xiebing%Bao%NULL%1,           Yongsheng%Li%dr_ysli@126.com%1,           Jun%Jin%jinjun0514@163.com%1]'
$ws.Cells.Item(26, 5).Value2 = '[Chenze%Li%NULL%0,           Jiangang%Jiang%NULL%0,           Feng%Wang%NULL%0,           Ning%Zhou%NULL%0,           Giacomo%Veronese%NULL%0,           Javid J.%Moslehi%NULL%0,           Enrico%Ammirati%NULL%0,           Dao Wen%Wang%NULL%0]'
$ws.Cells.Item(27, 5).Value2 = '[Miao%Luo%xref no email%0,    Jing%Liu%xref no email%0,    Weiling%Jiang%xref no email%0,    Shuang%Yue%xref no email%0,    Huiguo%Liu%xref no email%0,    Shuang%Wei%xref no email%0]'
$ws.Cells.Item(28, 5).Value2 = '[Zhao-Hua%Wang%NULL%1,           Chang%Shu%NULL%0,           Xiao%Ran%NULL%1,           Cui-Hong%Xie%NULL%1,           Lei%Zhang%NULL%0]'
$ws.Cells.Item(29, 5).Value2 = '[Benping%Zhang%NULL%1,           Chen%Dong%NULL%1,           Shengzhong%Li%NULL%1,           Xiaoqing%Song%NULL%1,           Wang%Wei%NULL%1,           Li%Liu%NULL%0]'
$ws.Cells.Item(30, 5).Value2 = '[Lin%Zhang%NULL%1,           Bin%Huang%NULL%1,           Hongzhen%Xia%NULL%1,           Hua%Fan%NULL%1,           Muxin%Zhu%NULL%1,           Liping%Zhu%NULL%1,           Huan%Zhang%NULL%1,           Xiaogen%Tao%NULL%1,           Shaohui%Cheng%NULL%1,           Jian%Chen%NULL%0]'
